# gtsummary_translated.xlsx -- add a new "Sum" translation row to Sheet1
# (adds digits= support for tbl_cross(); this row documents the
# translated strings for the new "Sum" statistic label).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at 70, pushing the existing rows 70:86 down to 71:87 ---
$ws.Rows("70:70").Insert()

# --- 2. Populate the new row 70. Shared strings must be registered in the
#        same order the final file uses them (774..781) so cell <v> indices
#        line up; that means writing the *first* occurrence of each new
#        string before any repeats. ---
$ws.Range("A70").Value2 = "tbl_summary"
$ws.Range("B70").Value2 = "Sum"
$ws.Range("K70").Value2 = "Summe"
$ws.Range("L70").Value2 = "Somme"
$ws.Range("M70").Value2 = "Soma"
$ws.Range("I70").Value2 = "Summa"
$ws.Range("O70").Value2 = "합집합"
$ws.Range("F70").Value2 = "和"
$ws.Range("J70").Value2 = "Suma"
$ws.Range("G70").Value2 = "和"
$ws.Range("H70").Value2 = "和"
$ws.Range("N70").Value2 = "Summa"

# C70/D70/E70 stay blank (no value), matching the other "Sum" columns.

# Match the wrapped-text formatting used by columns I/J elsewhere in the table.
$ws.Range("I70:J70").WrapText = $true

# Record the current selection the same way a user would after typing the row.
$ws.Range("J70").Select()

# --- 3. Extend the autofilter + its backing defined name over the new row ---
$ws.AutoFilterMode = $false
$ws.Range("A1:M87").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$M`$87"
    }
}
